$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.601.74"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = "'1.599.88"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "'212.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = "'0.514"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = "'26.86"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.87%  '
$ws.Range("D9").Value = "'0.251"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("E10").Value = '  +1.21%  '
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").Value = "'1.599.72"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.38%  '
$ws.Range("E14").Value = '  +3.05%  '
$ws.Range("D15").Value = "'29.607.26"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.60%  '
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = "'63.72"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.27%  '
$ws.Range("D18").Value = "'241.68"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.33%  '
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").Value = "'0.0₃0694"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("D24").Value = "'2.10"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.17%  '
$ws.Range("D25").Value = "'155.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("E26").Value = '  +1.52%  '
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("E30").Value = '  +2.20%  '
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("D32").Value = "'3.22"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("E33").Value = '  +3.01%  '
$ws.Range("D34").Value = "'1.423.32"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = "'1.54"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("D36").Value = "'2.87"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.68%  '
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("E39").Value = '  +3.05%  '
$ws.Range("E40").Value = '  +2.64%  '
$ws.Range("D41").Value = "'56.51"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +7.16%  '
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("E43").Value = '  +5.54%  '
$ws.Range("D44").Value = "'0.809"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").Value = "'0.991"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +16.79%  '
$ws.Range("D47").Value = "'66.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("D49").Value = "'1.740.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("D50").Value = "'86.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").Value = "'0.0₆0105"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.24%  '
